# customercare.xlsx - "New developed and ENW fixed scripts have been committed"
# Adds a new test case row (Customercare021) to the "Test Cases" sheet and
# tweaks a couple of cosmetic view/format settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Unicode punctuation used in the new description text (left single quote,
# right single quote) - build with char codes so encoding is unambiguous.
$lsquo = [char]0x2018
$rsquo = [char]0x2019

$tcid        = "Customercare021"
$jira        = "OPQA-5298 `n||OPQA-5299`n|| OPQA-5300"
$description = "Verify that Phone Number format should Only allow digits, " + $lsquo + "+" + $rsquo + " sign (at beginning), dashes " + $lsquo + "-" + $lsquo + ", parentheses '()' and spaces within the field`nVerify that spaces and special characters should be stripped out from the phone Number field before being submitted to salesforce`nVerify that Minimum 7 digits and Maximum of 40 characters should be inputted in the phone field in customer care page"

# New test case row 19: TCID | Jira id | Description | Runmode
$ws.Range("A19").Value = $tcid
$ws.Range("B19").Value = $jira
$ws.Range("C19").Value = $description
$ws.Range("D19").Value = "Y"

# Row grew to fit the wrapped, multi-line description.
$ws.Rows.Item(19).RowHeight = 45

# Column A widened slightly.
$ws.Columns.Item(1).ColumnWidth = 20

# View scrolled down a bit and selection moved as the author kept working.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C16").Select()
